$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.645.95"
Set-TextValue "E2" "  +1.12%  "

Set-TextValue "D3" "1.853.69"
Set-TextValue "E3" "  +0.74%  "

Set-TextValue "E4" "  -0.17%  "

Set-TextValue "D5" "264.55"
Set-TextValue "E5" "  +2.58%  "

Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  -0.19%  "

Set-TextValue "D7" "0.5222"
Set-TextValue "E7" "  +2.36%  "

Set-TextValue "D8" "0.3276"
Set-TextValue "E8" "  +1.15%  "

Set-TextValue "D9" "0.06797"
Set-TextValue "E9" "  +1.25%  "

Set-TextValue "D10" "18.86"
Set-TextValue "E10" "  -0.98%  "

Set-TextValue "D11" "0.7779"
Set-TextValue "E11" "  +1.33%  "

Set-TextValue "D12" "0.07764"
Set-TextValue "E12" "  +0.94%  "

Set-TextValue "D13" "1.856.93"
Set-TextValue "E13" "  +0.45%  "

Set-TextValue "E14" "  +0.86%  "

Set-TextValue "D15" "5.039"
Set-TextValue "E15" "  +0.42%  "

Set-TextValue "E16" "  -0.21%  "

Set-TextValue "D17" "14.03"
Set-TextValue "E17" "  -0.16%  "

Set-TextValue "B18" "Dai"
Set-TextValue "C18" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D18" "1.000"
Set-TextValue "E18" "  -0.08%  "

Set-TextValue "B19" "ShibaInu"
Set-TextValue "C19" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D19" "0.000007961"
Set-TextValue "E19" "  +1.30%  "

Set-TextValue "D20" "26.684.74"
Set-TextValue "E20" "  +1.01%  "

Set-TextValue "D21" "2.096.52"
Set-TextValue "E21" "  -1.68%  "

Set-TextValue "D22" "4.649"
Set-TextValue "E22" "  +1.69%  "

Set-TextValue "D23" "9.543"
Set-TextValue "E23" "  -0.41%  "

Set-TextValue "D24" "6.013"
Set-TextValue "E24" "  +0.74%  "

Set-TextValue "D25" "143.72"
Set-TextValue "E25" "  -1.03%  "

Set-TextValue "D26" "2.209"
Set-TextValue "E26" "  -5.18%  "

Set-TextValue "D27" "1.676"
Set-TextValue "E27" "  +1.54%  "

Set-TextValue "D28" "17.05"
Set-TextValue "E28" "  +0.50%  "

Set-TextValue "D29" "112.08"
Set-TextValue "E29" "  +1.15%  "

Set-TextValue "D30" "4.211"
Set-TextValue "E30" "  +0.09%  "

Set-TextValue "D31" "4.151"
Set-TextValue "E31" "  -0.38%  "

Set-TextValue "D32" "0.08765"
Set-TextValue "E32" "  +0.75%  "

Set-TextValue "D33" "0.04833"
Set-TextValue "E33" "  +0.54%  "

Set-TextValue "D34" "1.137"
Set-TextValue "E34" "  +0.62%  "

Set-TextValue "D35" "0.7195"
Set-TextValue "E35" "  +6.70%  "

Set-TextValue "D36" "2.873"
Set-TextValue "E36" "  +0.28%  "

Set-TextValue "D37" "3.115"
Set-TextValue "E37" "  +0.84%  "

Set-TextValue "D38" "0.01793"
Set-TextValue "E38" "  -0.81%  "

Set-TextValue "E39" "  -0.37%  "

Set-TextValue "D40" "0.4893"
Set-TextValue "E40" "  -0.33%  "

Set-TextValue "D41" "113.00"
Set-TextValue "E41" "  +0.40%  "

Set-TextValue "D42" "0.8973"
Set-TextValue "E42" "  -0.60%  "

Set-TextValue "D43" "6.087"
Set-TextValue "E43" "  -0.31%  "

Set-TextValue "D44" "1.001"
Set-TextValue "E44" "  -0.14%  "

Set-TextValue "D45" "7.732"
Set-TextValue "E45" "  -0.01%  "

Set-TextValue "D46" "0.4187"
Set-TextValue "E46" "  -1.32%  "

Set-TextValue "D47" "0.05926"
Set-TextValue "E47" "  +0.39%  "

Set-TextValue "D48" "9.079"
Set-TextValue "E48" "  -0.97%  "

Set-TextValue "B49" "Elrond"
Set-TextValue "C49" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D49" "35.11"
Set-TextValue "E49" "  -0.10%  "

Set-TextValue "B50" "Algorand"
Set-TextValue "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.1236"
Set-TextValue "E50" "  -3.50%  "

Set-TextValue "D51" "0.8870"
Set-TextValue "E51" "  +3.74%  "
